$wb = $excel.ActiveWorkbook

# ALC!row8: On the Drip / Eye Drops
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 16.6
$ws.Range("I8").Value = 16.6
$ws.Range("K8").Value = 49.8
$ws.Range("M8").Value = 89.19999999999999

# ALC!row32: Automata for the People / Crab Oil
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1460.8572
$ws.Range("I32").Value = 2266.6667
$ws.Range("J32").Value = 856.5
$ws.Range("K32").Value = 2266.6667
$ws.Range("L32").Value = 856.5
$ws.Range("M32").Value = -1940.6667
$ws.Range("N32").Value = -1508.5

# ALC!row101: Edge of the Arcane / Cunning Craftsman's Tea
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 16102.679
$ws.Range("I101").Value = 267.5
$ws.Range("J101").Value = 20421.363
$ws.Range("K101").Value = 802.5
$ws.Range("L101").Value = 61264.08900000001
$ws.Range("M101").Value = 819.5
$ws.Range("N101").Value = -64508.08900000001

# ALC!row112: Making Ends Meet / Superior Spiritbond Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1737902.6
$ws.Range("I112").Value = 1250
$ws.Range("J112").Value = 2085233.1
$ws.Range("K112").Value = 3750
$ws.Range("L112").Value = 6255699.300000001
$ws.Range("M112").Value = -2642
$ws.Range("N112").Value = -6257915.300000001

# ALC!row137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 19750836
$ws.Range("I137").Value = 4630747.5
$ws.Range("J137").Value = 56863780
$ws.Range("K137").Value = 13892242.5
$ws.Range("L137").Value = 170591340
$ws.Range("M137").Value = -13889692.5
$ws.Range("N137").Value = -170596440

# ARM!row2: Ain't Got No Ingots / Bronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20219.666
$ws.Range("I2").Value = 25655.285
$ws.Range("J2").Value = 1195
$ws.Range("K2").Value = 25655.285
$ws.Range("L2").Value = 1195
$ws.Range("M2").Value = -25542.285
$ws.Range("N2").Value = -1421

# ARM!row45: Hollow Hallmarks / Mythril Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2755.1875
$ws.Range("I45").Value = 1483.1666
$ws.Range("K45").Value = 1483.1666
$ws.Range("M45").Value = -1106.1666

# ARM!row63: Rivets Run through It / Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1630
$ws.Range("I63").Value = 1390
$ws.Range("J63").Value = 2350
$ws.Range("K63").Value = 1390
$ws.Range("L63").Value = 2350
$ws.Range("M63").Value = -704
$ws.Range("N63").Value = -3722

# ARM!row66: A Riveting Revival (L) / Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1630
$ws.Range("I66").Value = 1390
$ws.Range("J66").Value = 2350
$ws.Range("K66").Value = 6950
$ws.Range("L66").Value = 11750
$ws.Range("M66").Value = -3518
$ws.Range("N66").Value = -18614

# ARM!row80: A Squire to Inspire / Titanium Hoplon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 22633.166
$ws.Range("J80").Value = 22633.166
$ws.Range("L80").Value = 22633.166
$ws.Range("N80").Value = -24629.166

# ARM!row83: All's Fair in Highborn Assassination (L) / Titanium Hoplon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 22633.166
$ws.Range("J83").Value = 22633.166
$ws.Range("L83").Value = 67899.49800000001
$ws.Range("N83").Value = -77883.49800000001

# ARM!row116: No Scope / Titanbronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 20219.666
$ws.Range("I116").Value = 25655.285
$ws.Range("J116").Value = 1195
$ws.Range("K116").Value = 25655.285
$ws.Range("L116").Value = 1195
$ws.Range("M116").Value = -23361.285
$ws.Range("N116").Value = -5783

# BSM!row3: Hells Bells / Bronze Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20219.666
$ws.Range("I3").Value = 25655.285
$ws.Range("J3").Value = 1195
$ws.Range("K3").Value = 25655.285
$ws.Range("L3").Value = 1195
$ws.Range("M3").Value = -25541.285
$ws.Range("N3").Value = -1423

# BSM!row82: Spirituality Inspector / Titanium Lump Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27141.5
$ws.Range("J82").Value = 29876
$ws.Range("L82").Value = 29876
$ws.Range("N82").Value = -30642

# BSM!row85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 27141.5
$ws.Range("J85").Value = 29876
$ws.Range("L85").Value = 29876
$ws.Range("N85").Value = -32528

# CRP!row50: The Arsenal of Theocracy / Cobalt Halberd
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 11345.777
$ws.Range("J50").Value = 11345.777
$ws.Range("L50").Value = 11345.777
$ws.Range("N50").Value = -12595.777

# CRP!row58: You Do the Heavy Lifting / Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5062378.5
$ws.Range("I58").Value = 15121.143
$ws.Range("J58").Value = 22727780
$ws.Range("K58").Value = 15121.143
$ws.Range("L58").Value = 22727780
$ws.Range("M58").Value = -14918.143
$ws.Range("N58").Value = -22728186

# CRP!row59: Bow Down to Magic / Crab Bow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 16632.834
$ws.Range("J59").Value = 16632.834
$ws.Range("L59").Value = 16632.834
$ws.Range("N59").Value = -18922.834

# CRP!row136: Turali Quality / Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5062378.5
$ws.Range("I136").Value = 15121.143
$ws.Range("J136").Value = 22727780
$ws.Range("K136").Value = 45363.429
$ws.Range("L136").Value = 68183340
$ws.Range("M136").Value = -42813.429
$ws.Range("N136").Value = -68188440

# CUL!row36: Love's Crumpets Lost / Crumpet
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 54000
$ws.Range("I36").Value = 501
$ws.Range("J36").Value = 67374.75
$ws.Range("K36").Value = 1503
$ws.Range("L36").Value = 202124.25
$ws.Range("M36").Value = -1334
$ws.Range("N36").Value = -202462.25

# CUL!row56: Culture Club / Crowned Pie
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3002.5
$ws.Range("I56").Value = 3002.5
$ws.Range("K56").Value = 3002.5
$ws.Range("M56").Value = -2472.5

# CUL!row111: Soup for the Soldier / Broad Bean Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 250002020
$ws.Range("J111").Value = 4030
$ws.Range("L111").Value = 12090
$ws.Range("N111").Value = -18224

# CUL!row112: Sweet Tooth / Caramels
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 5752.467
$ws.Range("I112").Value = 2265.4
$ws.Range("J112").Value = 7496
$ws.Range("K112").Value = 6796.200000000001
$ws.Range("L112").Value = 22488
$ws.Range("M112").Value = -5688.200000000001
$ws.Range("N112").Value = -24704

# GSM!row102: Put the Metal to the Peddle / Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9301.733
$ws.Range("I102").Value = 9301.733
$ws.Range("K102").Value = 9301.733
$ws.Range("M102").Value = -7679.733

# GSM!row107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 377.4762
$ws.Range("I107").Value = 254.07143
$ws.Range("J107").Value = 624.2857
$ws.Range("K107").Value = 254.07143
$ws.Range("L107").Value = 624.2857
$ws.Range("M107").Value = 1665.92857
$ws.Range("N107").Value = -4464.2857

# GSM!row123: Workplace Workout / Ametrine Ring of Fending
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10046.571
$ws.Range("J123").Value = 10046.571
$ws.Range("L123").Value = 10046.571
$ws.Range("N123").Value = -14946.571

# GSM!row132: On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 19916844
$ws.Range("I132").Value = 20638132
$ws.Range("J132").Value = 18185758
$ws.Range("K132").Value = 61914396
$ws.Range("L132").Value = 54557274
$ws.Range("M132").Value = -61911866
$ws.Range("N132").Value = -54562334

# LTW!row22: Skin off Their Backs / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2589.3809
$ws.Range("I22").Value = 272.5
$ws.Range("J22").Value = 3134.5293
$ws.Range("K22").Value = 272.5
$ws.Range("L22").Value = 3134.5293
$ws.Range("M22").Value = 22.5
$ws.Range("N22").Value = -3724.5293

# LTW!row27: Fire and Hide / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2589.3809
$ws.Range("I27").Value = 272.5
$ws.Range("J27").Value = 3134.5293
$ws.Range("K27").Value = 272.5
$ws.Range("L27").Value = 3134.5293
$ws.Range("M27").Value = -165.5
$ws.Range("N27").Value = -3348.5293

# WVR!row113: A Tender Table / Pixie Floss
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 475.73685
$ws.Range("I113").Value = 499.91666
$ws.Range("J113").Value = 434.2857
$ws.Range("K113").Value = 1499.74998
$ws.Range("L113").Value = 1302.8571
$ws.Range("M113").Value = 670.2500199999999
$ws.Range("N113").Value = -5642.8571

# WVR!row132: Comfy Cabins / Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1835805.6
$ws.Range("I132").Value = 5685.5835
$ws.Range("J132").Value = 3404480
$ws.Range("K132").Value = 17056.7505
$ws.Range("L132").Value = 10213440
$ws.Range("M132").Value = -14526.7505
$ws.Range("N132").Value = -10218500
